# Fix typo 'Cyro-EM' to 'Cryo-EM' and update pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

# --- Update pav:createdOn value on the .metadata sheet ---
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2024-03-14T10:55:17-04:00"

# --- Fix the 'Cyro-EM' typo on the storage_medium lookup sheet ---
# Before:
#   A12/B12 = DMSO (serum) / https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125
#   A13/B13 = RNAlater     / http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348
#   A14/B14 = Cyro-EM      / https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333
# After (corrected spelling, moved up to keep alphabetic-ish ordering intact):
#   A12/B12 = Cryo-EM      / https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333
#   A13/B13 = DMSO (serum) / https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125
#   A14/B14 = RNAlater     / http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348
$storage = $wb.Worksheets.Item("storage_medium")

$storage.Range("A12").Value = "Cryo-EM"
$storage.Range("B12").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333"

$storage.Range("A13").Value = "DMSO (serum)"
$storage.Range("B13").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125"

$storage.Range("A14").Value = "RNAlater"
$storage.Range("B14").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348"
